# Applies the #5 commit: insurance sheet (保險) and investment sheet (事業投資)
# get the common trailing metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that every other sheet
# in this workbook already carries, plus a couple of header/label fixes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 5: 保險 (insurance)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Header row (row 1): B..E keep their meaning (company/name/owner) but the
# header labels themselves get corrected to the canonical field names, and
# F..K are brand-new metadata columns.
$ws5.Cells.Item(1, 2).Value = "company"
$ws5.Cells.Item(1, 3).Value = "name"
$ws5.Cells.Item(1, 4).Value = "owner"
$ws5.Cells.Item(1, 5).Value = "property_category"
$ws5.Cells.Item(1, 6).Value = "category"
$ws5.Cells.Item(1, 7).Value = "date"
$ws5.Cells.Item(1, 8).Value = "legislator_name"
$ws5.Cells.Item(1, 9).Value = "legislator_id"
$ws5.Cells.Item(1, 10).Value = "source_file"
$ws5.Cells.Item(1, 11).Value = "index"

# Copy the header formatting (bold + border, style index 1) onto the new
# header cells F1:K1 so they match B1:E1.
$ws5.Range("E1").Copy() | Out-Null
$ws5.Range("F1:K1").PasteSpecial(-4122) | Out-Null

# Data rows 2-5: fill in the new trailing columns; row index values mirror
# column A (the original per-row id).
$insuranceRows = @(
    @{ Row = 2; Index = 114; Name = "國泰人壽鍾愛終身壽險(分期繳納）" },
    @{ Row = 3; Index = 115; Name = "國泰住院醫療終身健康保險(分期缴納）" },
    @{ Row = 4; Index = 116; Name = "國泰福壽養老保險（分期繳納）" },
    @{ Row = 5; Index = 117; Name = "南山康寧終身壽險（分期繳納）" }
)

foreach ($r in $insuranceRows) {
    $row = $r.Row
    $ws5.Cells.Item($row, 3).Value = $r.Name
    $ws5.Cells.Item($row, 5).Value = "insurance"
    $ws5.Cells.Item($row, 6).Value = "normal"
    $ws5.Cells.Item($row, 7).Value = "2012-04-18"
    $ws5.Cells.Item($row, 8).Value = "林明溱"
    $ws5.Cells.Item($row, 9).Value = 1706
    $ws5.Cells.Item($row, 10).Value = "tmp80511"
    $ws5.Cells.Item($row, 11).Value = $r.Index

    $ws5.Range("E" + $row).Copy() | Out-Null
    $ws5.Range("F" + $row + ":K" + $row).PasteSpecial(-4122) | Out-Null
}

# Row 5's company was mis-copied as 國泰人壽 in the source data; it is 南山人壽.
$ws5.Cells.Item(5, 2).Value = "南山人壽"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet 6: 事業投資 (business investment)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Header row (row 1)
$ws6.Cells.Item(1, 2).Value = "owner"
$ws6.Cells.Item(1, 3).Value = "company"
$ws6.Cells.Item(1, 4).Value = "address"
$ws6.Cells.Item(1, 5).Value = "total"
$ws6.Cells.Item(1, 6).Value = "register_date"
$ws6.Cells.Item(1, 7).Value = "register_reason"
$ws6.Cells.Item(1, 8).Value = "property_category"
$ws6.Cells.Item(1, 9).Value = "category"
$ws6.Cells.Item(1, 10).Value = "date"
$ws6.Cells.Item(1, 11).Value = "legislator_name"
$ws6.Cells.Item(1, 12).Value = "legislator_id"
$ws6.Cells.Item(1, 13).Value = "source_file"
$ws6.Cells.Item(1, 14).Value = "index"

$ws6.Range("E1").Copy() | Out-Null
$ws6.Range("F1:N1").PasteSpecial(-4122) | Out-Null

# Data row 2: shift the old F/G values (address/date label) into their new
# homes and append the common trailing metadata.
$ws6.Cells.Item(2, 4).Value = "南投縣中山街202號4樓"
$ws6.Cells.Item(2, 6).Value = "89年10月05日"
$ws6.Cells.Item(2, 7).Value = "投資"
$ws6.Cells.Item(2, 8).Value = "investment"
$ws6.Cells.Item(2, 9).Value = "normal"
$ws6.Cells.Item(2, 10).Value = "2012-04-18"
$ws6.Cells.Item(2, 11).Value = "林明溱"
$ws6.Cells.Item(2, 12).Value = 1706
$ws6.Cells.Item(2, 13).Value = "tmp80511"
$ws6.Cells.Item(2, 14).Value = 130

$ws6.Range("E2").Copy() | Out-Null
$ws6.Range("F2:N2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
